$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)

$ws.Range("D2").Value = '36.349.96'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '1.933.58'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("E6").Value = '  -2.87%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.26'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.358'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0836'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("D12").Value = '2.215.94'
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.06'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -8.37%  '
$ws.Range("E14").Value = '  -6.67%  '
$ws.Range("E15").Value = '  -4.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.68%  '
$ws.Range("D17").Value = '1.939.51'
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").Value = '36.280.60'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  -2.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '225.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.46%  '
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("E24").Value = '  -7.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.05'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("E28").Value = '  -2.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.33%  '
$ws.Range("E30").Value = '  -2.13%  '
$ws.Range("E31").Value = '  -6.35%  '
$ws.Range("E32").Value = '  -7.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0617'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.82%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("B36").Value = 'THORChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.91'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.08%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.87%  '
$ws.Range("E38").Value = '  -4.00%  '
$ws.Range("E39").Value = '  +1.81%  '
$ws.Range("E40").Value = '  +0.44%  '
$ws.Range("E41").Value = '  -0.80%  '
$ws.Range("E42").Value = '  -2.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.14'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.52'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.98%  '
$ws.Range("D45").Value = '1.322.65'
$ws.Range("E45").Value = '  -2.72%  '
$ws.Range("E46").Value = '  -6.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.93'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.87%  '
$ws.Range("E49").Value = '  -0.06%  '
$ws.Range("D50").Value = '2.107.27'
$ws.Range("E50").Value = '  -2.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.47%  '
